# "update diagram and tasks"
# The task-assignment table is updated:
#  - Fix the participant name "Гарицай Г." -> "Грицай Г." (cell A3)
#  - Change that participant's assigned task from "Модель ЗУР (ракета)"
#    to "ПБУ, ПУ" (cell B4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Грицай Г."
$ws.Range("B4").Value = "ПБУ, ПУ"

# Leave the selection on the last-edited cell, as in the authored workbook.
$ws.Range("B4").Select() | Out-Null
